$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, D, E, F, H, I, K across rows 2-25 (row index 0 = row2 ... 23 = row25)
$colB = @(0.3674614980064348,0.3359955865645929,0.3169488502687443,0.3092554506760052,0.3079820794124544,0.3168448184519548,0.3565550016713814,0.4366222356240428,0.4968292398426399,0.5245289083198088,0.5350634907247525,0.5327926550979782,0.5253946831734027,0.5208691308151572,0.4950253117592354,0.4792510712152307,0.4702074018721021,0.4671503786114215,0.4809272340358746,0.527566412292515,0.5583122419030531,0.5418782262478885,0.4801693625250039,0.4147224710025341)
$colC = @(0.1346372942845164,0.1206964374743222,0.1122080978610995,0.1087665882746478,0.1081961771710098,0.11216161388964,0.1298154225431745,0.1650208782557172,0.1912740684740584,0.2033083033278729,0.2078790215188917,0.2068940227618157,0.2036840635829549,0.2017196575122,0.190489483239304,0.1836238972719286,0.1796835811860547,0.1783509190750294,0.1843538585956992,0.204626533255805,0.2179554705475084,0.2108341331743304,0.1840238219345451,0.1554309907689628)
$colD = @(0.06677200631255076,0.06609435825877341,0.06567230349913444,0.06549876890357709,0.06546985916762793,0.06566996945730708,0.06653956872931843,0.06819933090913466,0.06939370025285996,0.06993219330966838,0.07013545783409114,0.0700917094420106,0.06994892880299375,0.06986138804505515,0.06935841496847672,0.06904865079529543,0.06887002524356234,0.06880946581604519,0.0690816727439767,0.06999088425916256,0.07058133181671877,0.07026652974286662,0.06906674519893841,0.06775491684362578)
$colE = @(0.4060931560985637,0.3542110107288039,0.3224739158650181,0.3095682165934761,0.3074268207168132,0.3222997572284498,0.3881780573897942,0.5184239506268256,0.6149319374782181,0.6590511050089987,0.6757920739398457,0.6721850501930078,0.6604276972342262,0.653230489394474,0.6120532520446602,0.5868499337547206,0.5723739275611877,0.5674759995104637,0.5895307515949781,0.6638801718570875,0.712671296246981,0.686611377141304,0.5883187113625468,0.4830581212896021)
$colF = @(1.862171915237468,1.797652293064999,1.75875393034174,1.743081029428765,1.740489285127637,1.758541839760355,1.839775648506262,2.00485337260605,2.129796567972107,2.187462106484872,2.209419891174662,2.20468547152359,2.189266148174369,2.179837197392175,2.126044825968222,2.093258495743044,2.074478504083658,2.068133229986927,2.096740586249979,2.193791874659865,2.257927289824693,2.223631675065747,2.095166117392637,1.95956267409062)
$colH = @(0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429,0.07973214163530429)
$colI = @(1.315761505247679,1.279320014550251,1.257351274567412,1.248499934283075,1.24703624971238,1.257231494432844,1.30311142993385,1.396362469679786,1.466958043331005,1.499545009991337,1.511954190808552,1.509278554078094,1.500564526718051,1.495235982498414,1.464838025116492,1.446311801886011,1.435700425350149,1.432115186909144,1.448279344114781,1.503122159732513,1.539369063414696,1.519986036632432,1.447389695021698,1.370775727204389)
$colK = @(0.4976014749059061,0.4521574689581485,0.424592391364115,0.4134434422689708,0.4115972184579277,0.4244416935987658,0.4818617372616245,0.5971885363222498,0.6836555483528173,0.7233856272117123,0.7384885023147376,0.7352332334718028,0.724626982252687,0.7181379260845517,0.6810671472083811,0.6584275321428663,0.6454430153351609,0.641053037322223,0.660833701224476,0.7277407139966385,0.7718066299605084,0.7482565151437655,0.6597457747720341,0.5656901078249632)

for ($i = 0; $i -lt 24; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value2 = $colB[$i]
    $ws.Cells.Item($row, 3).Value2 = $colC[$i]
    $ws.Cells.Item($row, 4).Value2 = $colD[$i]
    $ws.Cells.Item($row, 5).Value2 = $colE[$i]
    $ws.Cells.Item($row, 6).Value2 = $colF[$i]
    $ws.Cells.Item($row, 8).Value2 = $colH[$i]
    $ws.Cells.Item($row, 9).Value2 = $colI[$i]
    $ws.Cells.Item($row, 11).Value2 = $colK[$i]
}

Write-Host "Updated pl_mw values for 380 kV case"
